$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Inscritos" (column E), "Pagos" (column F) and "Inscrições homologadas" (column H) values
$ws.Range("E2").Value = 21

$ws.Range("F5").Value = 10
$ws.Range("H5").Value = 10

$ws.Range("E6").Value = 41

$ws.Range("E7").Value = 21

$ws.Range("E9").Value = 15

$ws.Range("E12").Value = 19

$ws.Range("E16").Value = 241
$ws.Range("F16").Value = 68
$ws.Range("H16").Value = 68

$ws.Range("E18").Value = 70
